$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the "2027-04-25 Anzac Day" row (currently row 21) first,
# so the row number for the earlier insertion point (row 8) stays correct.
$ws.Rows.Item(22).Insert()
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "2027-04-26"
$ws.Cells.Item(22, 2).Value = "Anzac Day (additional)"

# Insert a new row after the "2026-04-25 Anzac Day" row (row 8).
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "2026-04-27"
$ws.Cells.Item(9, 2).Value = "Anzac Day (additional)"
